# Project DesignFirst - update rule R30 "From" threshold in sheet "Rules"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")
$ws.Range("C10").Value = 100
